$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: student ids for columns E..H become actual numeric ids (I4 keeps "id6" label)
$ws.Range("E4").Value = 4645669
$ws.Range("F4").Value = 4536908
$ws.Range("G4").Value = 4538420
$ws.Range("H4").Value = 4559118

# Tuesday (row 8) hours for week 3.1
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0

# Wednesday (row 9) hours for week 3.1
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 4
$ws.Range("H9").Value = 4

# Thursday (row 10) hours for week 3.1
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 4
$ws.Range("H10").Value = 8

# Friday (row 11) hours for week 3.1
$ws.Range("D11").Value = 8.5
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 8.5
$ws.Range("G11").Value = 8.5

# Monday (row 15) hours for week 3.2
$ws.Range("D15").Value = 6
$ws.Range("H15").Value = 4

# Restore selection to match author's final workbook state
$ws.Range("D16:D19").Select()

$wb.Save()
